$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$found = $ws.Cells.Find("lemmalist-greek")
if ($found -ne $null) {
    $found.EntireRow.Delete()
}
